# Refresh the cryptocurrency price/volume snapshot on the active sheet.
# (column D = Price, column E = Volume(1h) change — both stored as text
#  since the source feed renders them as preformatted strings.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Assigning a numeric-looking string normally gets coerced to a number
    # by Excel, which would strip formatting like trailing zeros. Forcing
    # the cell to Text first keeps it a literal string, matching the feed.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "29.041.34"
$ws.Range("E2").Value = "  -0.30%  "
# Row 3
$ws.Range("D3").Value = "1.832.28"
$ws.Range("E3").Value = "  -0.05%  "
# Row 4
Set-TextValue "D4" "0.9989"
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
Set-TextValue "D5" "241.72"
$ws.Range("E5").Value = "  +0.16%  "
# Row 6
Set-TextValue "D6" "0.6292"
$ws.Range("E6").Value = "  -4.37%  "
# Row 7
$ws.Range("E7").Value = "  +0.00%  "
# Row 8
Set-TextValue "D8" "45.34"
$ws.Range("E8").Value = "  +1.36%  "
# Row 9
Set-TextValue "D9" "0.07613"
$ws.Range("E9").Value = "  +2.83%  "
# Row 10
Set-TextValue "D10" "0.2925"
$ws.Range("E10").Value = "  +0.01%  "
# Row 11
Set-TextValue "D11" "22.89"
$ws.Range("E11").Value = "  +0.10%  "
# Row 12
Set-TextValue "D12" "0.07653"
$ws.Range("E12").Value = "  -1.16%  "
# Row 13
$ws.Range("D13").Value = "1.825.80"
$ws.Range("E13").Value = "  -0.64%  "
# Row 14
Set-TextValue "D14" "4.966"
$ws.Range("E14").Value = "  -0.41%  "
# Row 15
Set-TextValue "D15" "0.6670"
$ws.Range("E15").Value = "  +0.07%  "
# Row 16
Set-TextValue "D16" "82.51"
$ws.Range("E16").Value = "  -0.77%  "
# Row 17
Set-TextValue "D17" "0.000009392"
$ws.Range("E17").Value = "  +9.26%  "
# Row 18
Set-TextValue "D18" "5.991"
$ws.Range("E18").Value = "  -1.96%  "
# Row 19
$ws.Range("D19").Value = "28.938.87"
$ws.Range("E19").Value = "  -0.70%  "
# Row 20
Set-TextValue "D20" "225.91"
$ws.Range("E20").Value = "  -0.18%  "
# Row 21
Set-TextValue "D21" "12.35"
$ws.Range("E21").Value = "  -0.72%  "
# Row 22
Set-TextValue "D22" "0.9997"
$ws.Range("E22").Value = "  -0.17%  "
# Row 23
Set-TextValue "D23" "7.223"
$ws.Range("E23").Value = "  +1.68%  "
# Row 24
Set-TextValue "D24" "1.0000"
$ws.Range("E24").Value = "  +0.00%  "
# Row 25
Set-TextValue "D25" "160.34"
$ws.Range("E25").Value = "  -0.50%  "
# Row 26
Set-TextValue "D26" "8.438"
$ws.Range("E26").Value = "  -1.60%  "
# Row 27
Set-TextValue "D27" "0.1368"
$ws.Range("E27").Value = "  -2.44%  "
# Row 28
Set-TextValue "D28" "17.87"
$ws.Range("E28").Value = "  -0.56%  "
# Row 29
Set-TextValue "D29" "1.499"
$ws.Range("E29").Value = "  -0.62%  "
# Row 30
Set-TextValue "D30" "4.074"
$ws.Range("E30").Value = "  -0.60%  "
# Row 31
Set-TextValue "D31" "4.038"
$ws.Range("E31").Value = "  -0.08%  "
# Row 32
Set-TextValue "D32" "1.199"
$ws.Range("E32").Value = "  +0.89%  "
# Row 33
Set-TextValue "D33" "0.05216"
$ws.Range("E33").Value = "  -1.22%  "
# Row 34
Set-TextValue "D34" "1.851"
$ws.Range("E34").Value = "  -0.64%  "
# Row 35
Set-TextValue "D35" "1.152"
$ws.Range("E35").Value = "  +0.88%  "
# Row 36
Set-TextValue "D36" "0.7334"
$ws.Range("E36").Value = "  -0.52%  "
# Row 37
Set-TextValue "D37" "2.598"
$ws.Range("E37").Value = "  -2.10%  "
# Row 38
$ws.Range("D38").Value = "1.278.87"
$ws.Range("E38").Value = "  -1.68%  "
# Row 39
Set-TextValue "D39" "2.760"
$ws.Range("E39").Value = "  +0.69%  "
# Row 40
Set-TextValue "D40" "0.01790"
$ws.Range("E40").Value = "  +0.27%  "
# Row 41
Set-TextValue "D41" "6.552"
$ws.Range("E41").Value = "  +8.99%  "
# Row 42
Set-TextValue "D42" "0.8925"
$ws.Range("E42").Value = "  -2.21%  "
# Row 43
Set-TextValue "D43" "1.001"
$ws.Range("E43").Value = "  +0.10%  "
# Row 44
Set-TextValue "D44" "101.83"
$ws.Range("E44").Value = "  -0.31%  "
# Row 45
$ws.Range("D45").Value = "1.973.68"
$ws.Range("E45").Value = "  -0.66%  "
# Row 46
Set-TextValue "D46" "65.02"
$ws.Range("E46").Value = "  +2.05%  "
# Row 48
$ws.Range("E48").Value = "  -0.38%  "
# Row 49
Set-TextValue "D49" "0.3990"
$ws.Range("E49").Value = "  -0.37%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "8.845"
$ws.Range("E50").Value = "  +1.87%  "
# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.05759"
$ws.Range("E51").Value = "  -1.46%  "
